$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '27.674.85'
$ws.Range("D3").Value = '1.588.09'
$ws.Range("E3").Value = '  -2.52%  '
$ws.Range("E4").Value = '  +0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  -3.35%  '
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.22'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.49%  '
$ws.Range("E9").Value = '  -1.74%  '
$ws.Range("E10").Value = '  -2.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0867'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.49%  '
$ws.Range("D12").Value = '1.814.08'
$ws.Range("E12").Value = '  -2.50%  '
$ws.Range("D13").Value = '1.596.79'
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("E14").Value = '  -4.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.530'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.66%  '
$ws.Range("D16").Value = '27.650.52'
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.45'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.39%  '
$ws.Range("E18").Value = '  -4.23%  '
$ws.Range("D19").Value = '0.0₃0697'
$ws.Range("E19").Value = '  -3.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.74%  '
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("E22").Value = '  -4.88%  '
$ws.Range("E23").Value = '  -3.57%  '
$ws.Range("E24").Value = '  -3.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.07%  '
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.19%  '
$ws.Range("E30").Value = '  -2.44%  '
$ws.Range("E31").Value = '  -2.66%  '
$ws.Range("E32").Value = '  -5.34%  '
$ws.Range("D33").Value = '1.371.45'
$ws.Range("E33").Value = '  -3.26%  '
$ws.Range("E34").Value = '  -5.73%  '
$ws.Range("E35").Value = '  -4.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.975'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.03%  '
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("E39").Value = '  -3.24%  '
$ws.Range("E40").Value = '  -3.36%  '
$ws.Range("E41").Value = '  +0.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.972'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.07%  '
$ws.Range("E43").Value = '  -2.47%  '
$ws.Range("E44").Value = '  +2.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.23%  '
$ws.Range("D46").Value = '1.725.08'
$ws.Range("E46").Value = '  -2.55%  '
$ws.Range("E47").Value = '  -5.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("E49").Value = '  -1.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0967'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.46%  '
$ws.Range("E51").Value = '  -1.62%  '
